# Sample Project / Main.xlsx — "Project Sample Project is saved." commit.
#
# The only functional change in this commit is cell B11 on the "Rules"
# sheet: it used to hold the text "R40" and now holds the text "1"
# (still a text value, not a number — Excel stores it via the shared
# strings table, same as before). We reproduce that by writing the
# value with a leading apostrophe, which is how Excel/COM signals
# "store this as text even though it looks like a number" so the
# digit string doesn't silently get reinterpreted as the numeric value 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("B11").Value = "'1"
